$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task entries for row 7 (water dropping) and row 8 (art overhaul + ui, world beginnings)
$ws.Range("C7").Value = 0.3
$ws.Range("D7").Value = "water dropping"

$ws.Range("C8").Value = 5.0999999999999996
$ws.Range("D8").Value = "art overhaul + ui, world beginnings "

# Update selection to D10, matching the new active cell in the saved file
$ws.Range("D10").Select()
